# "Relatorio placas negativas" - fill in counts for the "305(-)" (column O)
# and "220(-)" (column R) plate-size buckets on rows 2-9 of Sheet1, then
# leave the selection on R6 (matching the author's last cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O2").Value = 1120
$ws.Range("R2").Value = 885

$ws.Range("O3").Value = 877
$ws.Range("R3").Value = 1500

$ws.Range("O4").Value = 73
$ws.Range("R4").Value = 1540

$ws.Range("O5").Value = 148
$ws.Range("R5").Value = 304

$ws.Range("O6").Value = 906

$ws.Range("O7").Value = 1120

$ws.Range("O8").Value = 1250

$ws.Range("O9").Value = 1120

# Matches the workbook's recorded selection after the edit (sheet1.xml: <selection activeCell="R6" sqref="R6"/>)
$ws.Range("R6").Select()
